# Update the Years header row on "Inputs & Income Statement" to shift the
# forecast window back 5 years (2024A..2029F -> 2019A..2024F).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inputs & Income Statement")

$ws.Range("C2").Value = "2019A"
$ws.Range("D2").Value = "2020F"
$ws.Range("E2").Value = 2021
$ws.Range("F2").Value = 2022
$ws.Range("G2").Value = 2023
$ws.Range("H2").Value = 2024

# Restore the view: scrolled so row 3 is at the top of the frozen pane, with
# the final selection left on F16.
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Range("F16").Select() | Out-Null
